$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(40, 8).Value = 2769.3125
$ws_ALC.Cells.Item(40, 9).Value = 2847.182
$ws_ALC.Cells.Item(40, 10).Value = 2598
$ws_ALC.Cells.Item(40, 11).Value = 2847.182
$ws_ALC.Cells.Item(40, 12).Value = 2598
$ws_ALC.Cells.Item(40, 13).Value = -2672.182
$ws_ALC.Cells.Item(40, 14).Value = -2948
$ws_ALC.Cells.Item(48, 8).Value = 250529.25
$ws_ALC.Cells.Item(48, 9).Value = 1000017
$ws_ALC.Cells.Item(48, 11).Value = 3000051
$ws_ALC.Cells.Item(48, 13).Value = -2999759
$ws_ALC.Cells.Item(56, 8).Value = 250529.25
$ws_ALC.Cells.Item(56, 9).Value = 1000017
$ws_ALC.Cells.Item(56, 11).Value = 3000051
$ws_ALC.Cells.Item(56, 13).Value = -2999517
$ws_ALC.Cells.Item(70, 8).Value = 2390.6206
$ws_ALC.Cells.Item(70, 9).Value = 963.2
$ws_ALC.Cells.Item(70, 10).Value = 3920
$ws_ALC.Cells.Item(70, 11).Value = 2889.6
$ws_ALC.Cells.Item(70, 12).Value = 11760
$ws_ALC.Cells.Item(70, 13).Value = -2619.6
$ws_ALC.Cells.Item(70, 14).Value = -12300
$ws_ALC.Cells.Item(73, 8).Value = 2390.6206
$ws_ALC.Cells.Item(73, 9).Value = 963.2
$ws_ALC.Cells.Item(73, 10).Value = 3920
$ws_ALC.Cells.Item(73, 11).Value = 2889.6
$ws_ALC.Cells.Item(73, 12).Value = 11760
$ws_ALC.Cells.Item(73, 13).Value = -1953.6
$ws_ALC.Cells.Item(73, 14).Value = -13632
$ws_ALC.Cells.Item(126, 8).Value = 11996.667
$ws_ALC.Cells.Item(126, 10).Value = 11996.667
$ws_ALC.Cells.Item(126, 12).Value = 11996.667
$ws_ALC.Cells.Item(126, 14).Value = -21876.667
$ws_ALC.Cells.Item(127, 8).Value = 1330.5667
$ws_ALC.Cells.Item(127, 9).Value = 645.2857
$ws_ALC.Cells.Item(127, 10).Value = 1539.1305
$ws_ALC.Cells.Item(127, 11).Value = 1935.8571
$ws_ALC.Cells.Item(127, 12).Value = 4617.3915
$ws_ALC.Cells.Item(127, 13).Value = 3024.1429
$ws_ALC.Cells.Item(127, 14).Value = -14537.3915
$ws_ALC.Cells.Item(129, 8).Value = 1035.3077
$ws_ALC.Cells.Item(129, 10).Value = 1224.4333
$ws_ALC.Cells.Item(129, 12).Value = 3673.2999
$ws_ALC.Cells.Item(129, 14).Value = -13673.2999
$ws_ALC.Cells.Item(138, 8).Value = 3600.3572
$ws_ALC.Cells.Item(138, 9).Value = 1825.1578
$ws_ALC.Cells.Item(138, 10).Value = 4511.946
$ws_ALC.Cells.Item(138, 11).Value = 5475.4734
$ws_ALC.Cells.Item(138, 12).Value = 13535.838
$ws_ALC.Cells.Item(138, 13).Value = -335.4733999999999
$ws_ALC.Cells.Item(138, 14).Value = -23815.838
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(32, 8).Value = 6976.521
$ws_ARM.Cells.Item(32, 9).Value = 4872.317
$ws_ARM.Cells.Item(32, 10).Value = 19301.143
$ws_ARM.Cells.Item(32, 11).Value = 4872.317
$ws_ARM.Cells.Item(32, 12).Value = 19301.143
$ws_ARM.Cells.Item(32, 13).Value = -4585.317
$ws_ARM.Cells.Item(32, 14).Value = -19875.143
$ws_ARM.Cells.Item(43, 8).Value = 19125.666
$ws_ARM.Cells.Item(43, 10).Value = 19125.666
$ws_ARM.Cells.Item(43, 12).Value = 19125.666
$ws_ARM.Cells.Item(43, 14).Value = -19751.666
$ws_ARM.Cells.Item(61, 8).Value = 993.9167
$ws_ARM.Cells.Item(61, 9).Value = 864.35297
$ws_ARM.Cells.Item(61, 10).Value = 1308.5714
$ws_ARM.Cells.Item(61, 11).Value = 864.35297
$ws_ARM.Cells.Item(61, 12).Value = 1308.5714
$ws_ARM.Cells.Item(61, 13).Value = -652.35297
$ws_ARM.Cells.Item(61, 14).Value = -1732.5714
$ws_ARM.Cells.Item(74, 8).Value = 1110.4
$ws_ARM.Cells.Item(74, 9).Value = 1140.125
$ws_ARM.Cells.Item(74, 10).Value = 1037.2307
$ws_ARM.Cells.Item(74, 11).Value = 1140.125
$ws_ARM.Cells.Item(74, 12).Value = 1037.2307
$ws_ARM.Cells.Item(74, 13).Value = -266.125
$ws_ARM.Cells.Item(74, 14).Value = -2785.2307
$ws_ARM.Cells.Item(77, 8).Value = 1110.4
$ws_ARM.Cells.Item(77, 9).Value = 1140.125
$ws_ARM.Cells.Item(77, 10).Value = 1037.2307
$ws_ARM.Cells.Item(77, 11).Value = 5700.625
$ws_ARM.Cells.Item(77, 12).Value = 5186.1535
$ws_ARM.Cells.Item(77, 13).Value = -1332.625
$ws_ARM.Cells.Item(77, 14).Value = -13922.1535
$ws_ARM.Cells.Item(136, 8).Value = 993.9167
$ws_ARM.Cells.Item(136, 9).Value = 864.35297
$ws_ARM.Cells.Item(136, 10).Value = 1308.5714
$ws_ARM.Cells.Item(136, 11).Value = 2593.05891
$ws_ARM.Cells.Item(136, 12).Value = 3925.7142
$ws_ARM.Cells.Item(136, 13).Value = -43.0589100000002
$ws_ARM.Cells.Item(136, 14).Value = -9025.7142
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(81, 8).Value = 12707.5
$ws_BSM.Cells.Item(81, 10).Value = 8094.2856
$ws_BSM.Cells.Item(81, 12).Value = 8094.2856
$ws_BSM.Cells.Item(81, 14).Value = -10216.2856
$ws_BSM.Cells.Item(84, 8).Value = 12707.5
$ws_BSM.Cells.Item(84, 10).Value = 8094.2856
$ws_BSM.Cells.Item(84, 12).Value = 24282.8568
$ws_BSM.Cells.Item(84, 14).Value = -34890.8568
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(6, 8).Value = 1127.5
$ws_CUL.Cells.Item(6, 9).Value = 84
$ws_CUL.Cells.Item(6, 10).Value = 2866.6667
$ws_CUL.Cells.Item(6, 11).Value = 252
$ws_CUL.Cells.Item(6, 12).Value = 8600.000100000001
$ws_CUL.Cells.Item(6, 13).Value = -139
$ws_CUL.Cells.Item(6, 14).Value = -8826.000100000001
$ws_CUL.Cells.Item(68, 8).Value = 2286714.8
$ws_CUL.Cells.Item(68, 9).Value = 2667500.8
$ws_CUL.Cells.Item(68, 10).Value = 2000
$ws_CUL.Cells.Item(68, 11).Value = 8002502.399999999
$ws_CUL.Cells.Item(68, 12).Value = 6000
$ws_CUL.Cells.Item(68, 13).Value = -8001691.399999999
$ws_CUL.Cells.Item(68, 14).Value = -7622
$ws_CUL.Cells.Item(71, 8).Value = 2286714.8
$ws_CUL.Cells.Item(71, 9).Value = 2667500.8
$ws_CUL.Cells.Item(71, 10).Value = 2000
$ws_CUL.Cells.Item(71, 11).Value = 24007507.2
$ws_CUL.Cells.Item(71, 12).Value = 18000
$ws_CUL.Cells.Item(71, 13).Value = -24003451.2
$ws_CUL.Cells.Item(71, 14).Value = -26112
$ws_CUL.Cells.Item(131, 8).Value = 5444.6816
$ws_CUL.Cells.Item(131, 9).Value = 247.54546
$ws_CUL.Cells.Item(131, 10).Value = 10641.818
$ws_CUL.Cells.Item(131, 11).Value = 742.6363799999999
$ws_CUL.Cells.Item(131, 12).Value = 31925.454
$ws_CUL.Cells.Item(131, 13).Value = 4297.36362
$ws_CUL.Cells.Item(131, 14).Value = -42005.454
$ws_CUL.Cells.Item(132, 8).Value = 1485.1428
$ws_CUL.Cells.Item(132, 9).Value = 1613.5714
$ws_CUL.Cells.Item(132, 10).Value = 1228.2858
$ws_CUL.Cells.Item(132, 11).Value = 14522.1426
$ws_CUL.Cells.Item(132, 12).Value = 11054.5722
$ws_CUL.Cells.Item(132, 13).Value = -11992.1426
$ws_CUL.Cells.Item(132, 14).Value = -16114.5722
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(122, 8).Value = 2422.558
$ws_GSM.Cells.Item(122, 9).Value = 2200.182
$ws_GSM.Cells.Item(122, 10).Value = 3156.4
$ws_GSM.Cells.Item(122, 11).Value = 6600.545999999999
$ws_GSM.Cells.Item(122, 12).Value = 9469.200000000001
$ws_GSM.Cells.Item(122, 13).Value = -4150.545999999999
$ws_GSM.Cells.Item(122, 14).Value = -14369.2
$ws_GSM.Cells.Item(124, 8).Value = 39500
$ws_GSM.Cells.Item(124, 10).Value = 39500
$ws_GSM.Cells.Item(124, 12).Value = 39500
$ws_GSM.Cells.Item(124, 14).Value = -49320
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(46, 8).Value = 65229.438
$ws_LTW.Cells.Item(46, 9).Value = 127833.875
$ws_LTW.Cells.Item(46, 10).Value = 2625
$ws_LTW.Cells.Item(46, 11).Value = 127833.875
$ws_LTW.Cells.Item(46, 12).Value = 2625
$ws_LTW.Cells.Item(46, 13).Value = -127645.875
$ws_LTW.Cells.Item(46, 14).Value = -3001
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(45, 8).Value = 18373
$ws_WVR.Cells.Item(45, 10).Value = 18373
$ws_WVR.Cells.Item(45, 12).Value = 18373
$ws_WVR.Cells.Item(45, 14).Value = -19355
$ws_WVR.Cells.Item(74, 8).Value = 7366.9
$ws_WVR.Cells.Item(74, 10).Value = 5900
$ws_WVR.Cells.Item(74, 12).Value = 5900
$ws_WVR.Cells.Item(74, 14).Value = -7772
$ws_WVR.Cells.Item(77, 8).Value = 7366.9
$ws_WVR.Cells.Item(77, 10).Value = 5900
$ws_WVR.Cells.Item(77, 12).Value = 17700
$ws_WVR.Cells.Item(77, 14).Value = -27060